$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-looking string (DD.MM.YYYY) into a cell without
# letting the engine auto-convert it into a real date serial/value. We
# temporarily force a text number-format, assign the literal string, then
# clear the formatting again so the cell ends up with no explicit style
# (matching the plain text cells used elsewhere in this sheet).
function Set-TextDate($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Helper: write the "date added" (M column) value, which keeps the
# existing custom date/time number format (style index 2 in the original
# workbook, numFmtId 165 "YYYY-MM-DD HH:MM:SS").
function Set-AddedDate($addr, $serial) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $rng.Value = $serial
}

### Row 2 ###
$ws.Range("B2").Value = 1973
$ws.Range("C2").Value = "ПБ"
$ws.Range("D2").Value = "Использование библиотеки Facerecognition и фреймворка Django для распознавания лиц в реальном времени"
Set-TextDate "E2" "11.08.2023"
$ws.Range("F2").Value = "ББ"
$ws.Range("G2").Value = "ВВБ"
$ws.Range("H2").Value = "Белова Елена Витальевна"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "Использование библиотеки Facerecognition и фреймворка Django для распознавания лиц в реальном времени"
$ws.Range("K2").Value = 1973
$ws.Range("L2").Value = "Использование библиотеки facerecognition и фреймворка Django для распознавания лиц в реальном времени."
Set-AddedDate "M2" 45110

### Row 3 ###
$ws.Range("B3").Value = 1976
$ws.Range("C3").Value = "СибБ"
$ws.Range("D3").Value = "Microsoft Bing и ChatGPT - buddy DA-DS-аудитора"
Set-TextDate "E3" "15.08.2023"
$ws.Range("F3").Value = "СЗБ"
$ws.Range("G3").Value = "ББ"
$ws.Range("H3").Value = "Шайдурова Арина Владимировна"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = "Microsoft Bing и ChatGPT - buddy DA-DS-аудитора"
$ws.Range("K3").Value = 1976
$ws.Range("L3").Value = "Microsoft Bing и ChatGPT - buddy будущего"
Set-AddedDate "M3" 45030

### Row 4 ###
$ws.Range("B4").Value = 1974
$ws.Range("C4").Value = "СРБ"
$ws.Range("D4").Value = "Как скачать отчёты из Power BI на диск с помощью PowerShell"
Set-TextDate "E4" "17.08.2023"
$ws.Range("F4").Value = "СЗБ"
$ws.Range("G4").Value = "ЦЧБ"
$ws.Range("H4").Value = "Ермолаева Светлана Алексеевна"
# I4 stays 1 (unchanged)
$ws.Range("J4").Value = "Как скачать отчёты из Power BI на диск с помощью PowerShell"
$ws.Range("K4").Value = 1974
$ws.Range("L4").Value = "Как сохранить все отчёты с сервера Power BI"
Set-AddedDate "M4" 45079

### Row 5 ###
$ws.Range("A5").Value = "Habr"
$ws.Range("B5").Value = 1965
$ws.Range("C5").Value = "ЦЧБ"
$ws.Range("D5").Value = "Замена Paint в задачах разметки графических данных"
Set-TextDate "E5" "14.08.2023"
$ws.Range("F5").Value = "ЮЗБ"
$ws.Range("G5").Value = "МБ"
$ws.Range("H5").Value = "Ермолаева Светлана Алексеевна"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = "Замена Paint в задачах разметки графических данных"
$ws.Range("K5").Value = 1965
$ws.Range("L5").Value = "Замена paint в задачах разметки графических данных"
Set-AddedDate "M5" 45062
